$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.148.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.22%  '

# Row 3 - Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.84%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.43%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.43%  '

# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4640'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.80%  '

# Row 8 - Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3709'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.22%  '

# Row 9 - Dogecoin
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07379'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.56%  '

# Row 10 - Polygon
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8824'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.32%  '

# Row 11 - TRON
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07915'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.86%  '

# Row 12 - Solana
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.93'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.23%  '

# Row 13 - WrappedEther
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.858.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.30%  '

# Row 14 - Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.370'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.29%  '

# Row 15 - Chainlink
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.595'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.86%  '

# Row 16 - Litecoin
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.85%  '

# Row 17 - BinanceUSD
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.005'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.42%  '

# Row 18 - ShibaInu
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008941'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.51%  '

# Row 19 - Dai
$ws.Range('E19').Value = '  -0.44%  '

# Row 20 - Avalanche
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.09%  '

# Row 21 - WrappedBTC
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.179.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.84%  '

# Row 22 - Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.136'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '

# Row 23 - Cosmos
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.94%  '

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.075.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.59%  '

# Row 25 - Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.57%  '

# Row 26 - Toncoin
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.869'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.37%  '

# Row 27 - EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.02%  '

# Row 28 - LidoDAOToken
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.075'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.71%  '

# Row 29 - InternetComputer(DFINITY)->BitcoinCash
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.65'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.79%  '

# Row 30 - BitcoinCash->InternetComputer(DFINITY)
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.124'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.30%  '

# Row 31 - Stellar
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08881'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.29%  '

# Row 32 - HuobiToken
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.962'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.28%  '

# Row 33 - ImmutableX
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7386'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.95%  '

# Row 34 - Filecoin
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.465'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '

# Row 35 - ARBITRUM
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.51%  '

# Row 36 - TrustWalletToken->RenderToken
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.523'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.48%  '

# Row 37 - RenderToken->TrustWalletToken
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.079'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.61%  '

# Row 38 - VeChain
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01951'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.50%  '

# Row 39 - Hedera
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05261'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.75%  '

# Row 40 - MXToken
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.970'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.27%  '

# Row 41 - FraxShare
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.068'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.49%  '

# Row 42 - TheSandbox
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5164'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.05%  '

# Row 43 - Algorand
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1634'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.30%  '

# Row 44 - Aptos
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.189'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.95%  '

# Row 45 - Decentraland
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4849'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.91%  '

# Row 46 - PaxDollar->EnergySwap
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.99%  '

# Row 47 - EnergySwap->PaxDollar
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.60%  '

# Row 48 - Quant
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.02%  '

# Row 49 - NEARProtocol
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.627'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.19%  '

# Row 50 - Cronos
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06229'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.58%  '

# Row 51 - Aave
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.36%  '
